$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: updated company data ---
$ws.Cells.Item(2, 1).Value = 'Tunisia'
$ws.Cells.Item(2, 2).Value = "'2"
$ws.Cells.Item(2, 3).Value = 'Insurance (General)'
$ws.Cells.Item(2, 4).Value = 0.08165
$ws.Cells.Item(2, 5).Value = -0.043
$ws.Cells.Item(2, 7).Value = 0.1614035087719298
$ws.Cells.Item(2, 8).Value = 0.1614035087719298
$ws.Cells.Item(2, 9).Value = 0.1092731829573935
$ws.Cells.Item(2, 10).Value = 0.07221157941196968
$ws.Cells.Item(2, 11).Value = 13.59
$ws.Cells.Item(2, 12).Value = 0.06812030075187969
$ws.Cells.Item(2, 13).Value = 0
$ws.Cells.Item(2, 14).Value = 0
$ws.Cells.Item(2, 15).Value = 0
$ws.Cells.Item(2, 16).Value = 0
$ws.Cells.Item(2, 17).Value = 0
$ws.Cells.Item(2, 18).Value = 0
$ws.Cells.Item(2, 19).Value = 0
$ws.Cells.Item(2, 21).Value = 14.18
$ws.Cells.Item(2, 22).Value = 0.05962994112699747
$ws.Cells.Item(2, 23).Value = 0.0914950980392157
$ws.Cells.Item(2, 24).Value = 0.07572266396924258
$ws.Cells.Item(2, 25).Value = 0.01577243406997311
$ws.Cells.Item(2, 26).Value = 1.263738003990752
$ws.Cells.Item(2, 27).Value = 0.1043391827583219
$ws.Cells.Item(2, 28).Value = 0.07572266396924258
$ws.Cells.Item(2, 29).Value = 0.0286165187890793
$ws.Cells.Item(2, 30).Value = 0
$ws.Cells.Item(2, 31).Value = 0
$ws.Cells.Item(2, 32).Value = 0
$ws.Cells.Item(2, 33).Value = -14.18
$ws.Cells.Item(2, 34).Value = 0
$ws.Cells.Item(2, 35).Value = 0
$ws.Cells.Item(2, 36).Value = -0.06341114390483857
$ws.Cells.Item(2, 37).Value = -0.0891151332327803
$ws.Cells.Item(2, 38).Value = 0
$ws.Cells.Item(2, 39).Value = 0
$ws.Cells.Item(2, 40).Value = 0
$ws.Cells.Item(2, 42).Value = -0.6274336283185841

# --- Row 3: updated company data ---
$ws.Cells.Item(3, 1).Value = 'Tunisia'
$ws.Cells.Item(3, 2).Value = 'Compagnie d''Assurances et de Réassurances ASTREE (BVMT:AST)'
$ws.Cells.Item(3, 3).Value = 'Insurance (General)'
$ws.Cells.Item(3, 4).Value = 0.138
$ws.Cells.Item(3, 5).Value = 0.107
$ws.Cells.Item(3, 7).Value = 0.3567839195979899
$ws.Cells.Item(3, 8).Value = 0.3567839195979899
$ws.Cells.Item(3, 9).Value = 0.1926298157453936
$ws.Cells.Item(3, 10).Value = 0.122470148909686
$ws.Cells.Item(3, 11).Value = 6.94
$ws.Cells.Item(3, 12).Value = 0.1162479061976549
$ws.Cells.Item(3, 13).Value = -0
$ws.Cells.Item(3, 14).Value = -0
$ws.Cells.Item(3, 15).Value = -0
$ws.Cells.Item(3, 16).Value = -0
$ws.Cells.Item(3, 17).Value = -0
$ws.Cells.Item(3, 18).Value = -0
$ws.Cells.Item(3, 19).Value = 0
$ws.Cells.Item(3, 21).Value = 3.88
$ws.Cells.Item(3, 22).Value = 0.02928301886792453
$ws.Cells.Item(3, 23).Value = 0.1275735294117647
$ws.Cells.Item(3, 24).Value = 0.07572266396924258
$ws.Cells.Item(3, 25).Value = 0.05185086544252214
$ws.Cells.Item(3, 26).Value = 1.162270028229339
$ws.Cells.Item(3, 27).Value = 0.1423433834305121
$ws.Cells.Item(3, 28).Value = 0.07572266396924258
$ws.Cells.Item(3, 29).Value = 0.06662071946126956
$ws.Cells.Item(3, 30).Value = 0
$ws.Cells.Item(3, 31).Value = 0
$ws.Cells.Item(3, 32).Value = 0
$ws.Cells.Item(3, 33).Value = -3.88
$ws.Cells.Item(3, 34).Value = 0
$ws.Cells.Item(3, 35).Value = 0
$ws.Cells.Item(3, 36).Value = -0.03016638158917742
$ws.Cells.Item(3, 37).Value = -0.08580274214949138
$ws.Cells.Item(3, 38).Value = 0
$ws.Cells.Item(3, 39).Value = 0
$ws.Cells.Item(3, 40).Value = 0
$ws.Cells.Item(3, 42).Value = -0.3316239316239316

# --- Row 4: updated company data ---
$ws.Cells.Item(4, 1).Value = 'Tunisia'
$ws.Cells.Item(4, 2).Value = 'Société Tunisienne d''Assurances et de Réassurances (BVMT:STAR)'
$ws.Cells.Item(4, 3).Value = 'Insurance (General)'
$ws.Cells.Item(4, 4).Value = 0.0253
$ws.Cells.Item(4, 5).Value = -0.193
$ws.Cells.Item(4, 7).Value = 0.07796852646638054
$ws.Cells.Item(4, 8).Value = 0.07796852646638054
$ws.Cells.Item(4, 9).Value = 0.0736766809728183
$ws.Cells.Item(4, 10).Value = 0.05053416017291143
$ws.Cells.Item(4, 11).Value = 6.65
$ws.Cells.Item(4, 12).Value = 0.04756795422031473
$ws.Cells.Item(4, 13).Value = -0
$ws.Cells.Item(4, 14).Value = -0
$ws.Cells.Item(4, 15).Value = -0
$ws.Cells.Item(4, 16).Value = -0
$ws.Cells.Item(4, 17).Value = -0
$ws.Cells.Item(4, 18).Value = -0
$ws.Cells.Item(4, 19).Value = 0
$ws.Cells.Item(4, 21).Value = 10.3
$ws.Cells.Item(4, 22).Value = 0.09781576448243116
$ws.Cells.Item(4, 23).Value = 0.05541666666666667
$ws.Cells.Item(4, 24).Value = 0.07572266396924258
$ws.Cells.Item(4, 25).Value = -0.02030599730257591
$ws.Cells.Item(4, 26).Value = 1.312676056338028
$ws.Cells.Item(4, 27).Value = 0.06633498208613163
$ws.Cells.Item(4, 28).Value = 0.07572266396924258
$ws.Cells.Item(4, 29).Value = -0.009387681883110954
$ws.Cells.Item(4, 30).Value = 0
$ws.Cells.Item(4, 31).Value = 0
$ws.Cells.Item(4, 32).Value = 0
$ws.Cells.Item(4, 33).Value = -10.3
$ws.Cells.Item(4, 34).Value = 0
$ws.Cells.Item(4, 35).Value = 0
$ws.Cells.Item(4, 36).Value = -0.108421052631579
$ws.Cells.Item(4, 37).Value = -0.09043020193151888
$ws.Cells.Item(4, 38).Value = 0
$ws.Cells.Item(4, 39).Value = 0
$ws.Cells.Item(4, 40).Value = 0
$ws.Cells.Item(4, 42).Value = -0.944954128440367

# buybacks_cash_returned (column T) is dropped from the refreshed dataset
$ws.Range("T2:T4").ClearContents()

# Old STAR row (row 5) data has been merged into row 4; drop the now-duplicate row
$ws.Rows(5).Delete()
